# Apply the bimonthly "directed trips" test-run adjustments described in the
# commit "adjustments and test runs":
#   - F14:F29, F41:F55, F109:F124  : 16.5 -> 17
#   - F56:F64 and R56:R64          : 16.5 -> 15
#   - Selection left on F109:F124 (the last block touched), matching the
#     sheet's saved cursor/selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Value adjustments -----------------------------------------------------
$ws.Range("F14:F29").Value = 17
$ws.Range("F41:F55").Value = 17
$ws.Range("F56:F64").Value = 15
$ws.Range("R56:R64").Value = 15
$ws.Range("F109:F124").Value = 17

# --- Leave the selection / active cell where the author left it ------------
$ws.Range("F109:F124").Select()
$excel.ActiveWindow.ScrollRow = 103
$excel.ActiveWindow.ScrollColumn = 1
